$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("C3").Value = -0.06778564228962507
$ws.Range("D3").Value = -0.6473335752170566
$ws.Range("E3").Value = 0.5493948348744779
$ws.Range("F3").Value = 0.06853313509662716
$ws.Range("G3").Value = -0.4539519709241162
$ws.Range("H3").Value = 0.4222925583339146
$ws.Range("I3").Value = -0.01318125193054271
$ws.Range("J3").Value = 0.9636801478230531
$ws.Range("K3").Value = 0.02831260680176752
$ws.Range("L3").Value = 0.06040454062648083
$ws.Range("M3").Value = -0.4947404136212541
$ws.Range("N3").Value = 0.4341926614871541
$ws.Range("O3").Value = -0.06040454658422392
$ws.Range("P3").Value = -0.4947401954287738
$ws.Range("Q3").Value = 0.5658073389676668
$ws.Range("R3").Value = 19.20749420500179
$ws.Range("S3").Value = -0.8909745524311108
$ws.Range("T3").Value = 270.7502469379378
$ws.Range("U3").Value = -0.2706257499234986
$ws.Range("V3").Value = 278.3452178247458
$ws.Range("Z3").Value = 1.288391225169693
$ws.Range("AA3").Value = 40.80399446376837
$ws.Range("AE3").Value = 1.187869740427624
$ws.Range("AF3").Value = 23.09551280799815
$ws.Range("AJ3").Value = 0.6970048778360557
$ws.Range("AO3").Value = 0.004284105657490979
$ws.Range("AP3").Value = 0.1409865461848463
$ws.Range("AT3").Value = 0.003949855739399389
$ws.Range("AU3").Value = 0.079799946695386
$ws.Range("AY3").Value = 0.002317652031542643
$ws.Range("BA3").Value = 0.9848193982062357
$ws.Range("BB3").Value = 0.7328930431439958
$ws.Range("BC3").Value = 12.44645988964466
$ws.Range("BD3").Value = -0.02670952659263602
$ws.Range("BE3").Value = 13.23220338983051
$ws.Range("BF3").Value = 1.619311068683604
$ws.Range("BG3").Value = -0.2760595240799184
$ws.Range("BH3").Value = 5.465107967864329
$ws.Range("BI3").Value = -0.06189362945645832
$ws.Range("BJ3").Value = 7.252542372881357
$ws.Range("BK3").Value = -1.457444750963728
$ws.Range("BL3").Value = -1.044297698570743
$ws.Range("BM3").Value = 6.554914141061038
$ws.Range("BN3").Value = 0.04091174751607241
$ws.Range("BO3").Value = 5.389830508474576
$ws.Range("BP3").Value = -0.1852502881121705
$ws.Range("BQ3").Value = 1.040731765825411
$ws.Range("BR3").Value = 0.4264377837557281
$ws.Range("BS3").Value = -0.005727644652250154
$ws.Range("BT3").Value = 0.5898305084745763
$ws.Range("BU3").Value = 0.07013598475023911
$ws.Range("BV3").Value = -0.439494650126642
$ws.Range("BW3").Value = 0.4233399886896293
$ws.Range("BX3").Value = -0.002816439246258236
$ws.Range("BY3").Value = 0.5045109365732364
$ws.Range("BZ3").Value = -0.07172491627521779
$ws.Range("CA3").Value = -0.6628430903675414
$ws.Range("CB3").Value = 0.5423709256652932
$ws.Range("CC3").Value = 0.003235143316032676
$ws.Range("CD3").Value = 0.4492998489837577
$ws.Range("CE3").Value = -0.01588796965796603
$ws.Range("CF3").Value = 1.119177552663611
$ws.Range("CG3").Value = 0.03428908589053887
$ws.Range("CH3").Value = -0.0004187040697744526
$ws.Range("CI3").Value = 0.04618921444300571
$ws.Range("CJ3").Value = 0.1764054364261761
$ws.Range("CK3").Value = -0.5700887690516562
$ws.Range("CL3").Value = 3.052461139162825
$ws.Range("CM3").Value = -0.008766803039158428
$ws.Range("CN3").Value = 3.305084745762713

# Row 10
$ws.Range("C10").Value = 0.04594943933365733
$ws.Range("D10").Value = 1.074912531447437
$ws.Range("E10").Value = 0.4564152473620577
$ws.Range("F10").Value = 0.0244044040211456
$ws.Range("G10").Value = -0.008021278785228942
$ws.Range("H10").Value = 0.448111535110376
$ws.Range("I10").Value = -0.06130373665863323
$ws.Range("J10").Value = 0.7155785127028472
$ws.Range("K10").Value = 0.09547320969296454
$ws.Range("L10").Value = 0.0256641040484506
$ws.Range("M10").Value = -0.9680906120826636
$ws.Range("N10").Value = 0.4874862825823784
$ws.Range("O10").Value = -0.02566410606166351
$ws.Range("P10").Value = -0.9680904241389079
$ws.Range("Q10").Value = 0.5125137174503233
$ws.Range("R10").Value = 19.68630553591419
$ws.Range("S10").Value = 1.242702135667751
$ws.Range("T10").Value = 279.4974579697418
$ws.Range("U10").Value = 1.271501404109685
$ws.Range("V10").Value = 242.9397850870029
$ws.Range("W10").Value = 26.50902810630701
$ws.Range("X10").Value = 1.203454915746716
$ws.Range("Y10").Value = 276.8657744925675
$ws.Range("Z10").Value = 1.327268597059798
$ws.Range("AA10").Value = 238.7943155145231
$ws.Range("AE10").Value = 1.289604193021523
$ws.Range("AF10").Value = 34.64489657338449
$ws.Range("AK10").Value = 36.72629178769042
$ws.Range("AL10").Value = -0.01381004801219559
$ws.Range("AM10").Value = -1.545264415494619
$ws.Range("AN10").Value = 0.9851059480952336
$ws.Range("AO10").Value = 0.00032775304165062
$ws.Range("AP10").Value = 0.9758351022102892
$ws.Range("AT10").Value = 0.004420710677738483
$ws.Range("AU10").Value = 0.1457569661642358
$ws.Range("AZ10").Value = 0.1545137494666013
$ws.Range("BA10").Value = 0.5731573103274017
$ws.Range("BB10").Value = -0.2976183881606958
$ws.Range("BC10").Value = 15.47742500097059
$ws.Range("BD10").Value = -0.01788427819988324
$ws.Range("BE10").Value = 15.99322033898305
$ws.Range("BF10").Value = 0.7610745062273941
$ws.Range("BG10").Value = -0.7973505130943954
$ws.Range("BH10").Value = 6.89058567184697
$ws.Range("BI10").Value = -0.009234365867913617
$ws.Range("BJ10").Value = 7.149152542372884
$ws.Range("BK10").Value = 0.8942249191081975
$ws.Range("BL10").Value = 0.971624588937006
$ws.Range("BM10").Value = 7.190873666392797
$ws.Range("BN10").Value = -0.0006428988895383258
$ws.Range("BO10").Value = 7.222033898305084
$ws.Range("BP10").Value = -0.6076219650536202
$ws.Range("BQ10").Value = 0.6381644647234
$ws.Range("BR10").Value = 1.395965660347877
$ws.Range("BS10").Value = -0.008007013442431338
$ws.Range("BT10").Value = 1.622033898305085
$ws.Range("BU10").Value = 0.02305447746770857
$ws.Range("BV10").Value = -0.533297627246399
$ws.Range("BW10").Value = 0.4407721488581078
$ws.Range("BX10").Value = 0.0001151750045676211
$ws.Range("BY10").Value = 0.4372334239169262
$ws.Range("BZ10").Value = 0.04710129651446118
$ws.Range("CA10").Value = 1.178400101372425
$ws.Range("CB10").Value = 0.4611553558822211
$ws.Range("CC10").Value = 0.0005088899298941564
$ws.Range("CD10").Value = 0.447135198436623
$ws.Range("CE10").Value = -0.04944739260764445
$ws.Range("CF10").Value = 0.698577242014558
$ws.Range("CG10").Value = 0.09807249550563846
$ws.Range("CH10").Value = -0.0006240649344617842
$ws.Range("CI10").Value = 0.1156313776464508
$ws.Range("CJ10").Value = 0.6397601563604972
$ws.Range("CK10").Value = -0.1605026579999047
$ws.Range("CL10").Value = 9.442410896626862
$ws.Range("CM10").Value = -0.02635885447106967
$ws.Range("CN10").Value = 10.20508474576271

# Row 47
$ws.Range("C47").Value = 0.04933434531231579
$ws.Range("D47").Value = -0.4566932577931364
$ws.Range("E47").Value = 0.3757131489031663
$ws.Range("F47").Value = 0.08548033307453555
$ws.Range("G47").Value = 0.1967785183102628
$ws.Range("H47").Value = 0.3716280170141263
$ws.Range("I47").Value = -0.1282082648105856
$ws.Range("J47").Value = -0.03934652836813769
$ws.Range("K47").Value = 0.2526588339328921
$ws.Range("L47").Value = 0.06395054454863799
$ws.Range("M47").Value = 0.6646357150833417
$ws.Range("N47").Value = 0.4837423983560422
$ws.Range("O47").Value = -0.06395054258265551
$ws.Range("P47").Value = 0.6646356251562832
$ws.Range("Q47").Value = 0.5162576014286395
$ws.Range("R47").Value = -141.7996975441891
$ws.Range("S47").Value = 0.3341659597316323
$ws.Range("T47").Value = 332.8563184753812
$ws.Range("U47").Value = 6.864197113260457
$ws.Range("V47").Value = 143.2604034068777
$ws.Range("Z47").Value = 0.8164895236040338
$ws.Range("AA47").Value = 28.2762393943345
$ws.Range("AE47").Value = 0.7660386939919104
$ws.Range("AF47").Value = 22.99560268954818
$ws.Range("AK47").Value = 26.39204763885191
$ws.Range("AO47").Value = 0.004298644131702872
$ws.Range("AP47").Value = 0.1483019548653907
$ws.Range("AT47").Value = 0.004033031216433099
$ws.Range("AU47").Value = 0.1206063078123157
$ws.Range("AZ47").Value = 0.1384198302736988
$ws.Range("BA47").Value = 2.849787013906935
$ws.Range("BB47").Value = -0.07869550531905882
$ws.Range("BC47").Value = 17.95042919101711
$ws.Range("BD47").Value = -0.1960013670539986
$ws.Range("BE47").Value = 23.33646616541353
$ws.Range("BF47").Value = 2.107738800822981
$ws.Range("BG47").Value = 0.09954330516220268
$ws.Range("BH47").Value = 6.996259553238923
$ws.Range("BI47").Value = -0.1313738892686262
$ws.Range("BJ47").Value = 10.61278195488722
$ws.Range("BK47").Value = 1.86810247399619
$ws.Range("BL47").Value = -0.2478590721108347
$ws.Range("BM47").Value = 7.204612495589292
$ws.Range("BN47").Value = -0.1052973342447028
$ws.Range("BO47").Value = 10.0921052631579
$ws.Range("BP47").Value = -1.059150562442817
$ws.Range("BQ47").Value = 0.008282737067073408
$ws.Range("BR47").Value = 3.732299514537933
$ws.Range("BS47").Value = 0.03885850991114145
$ws.Range("BT47").Value = 2.663533834586467
$ws.Range("BU47").Value = 0.07101354288319581
$ws.Range("BV47").Value = 0.2670116941100222
$ws.Range("BW47").Value = 0.3742732619501933
$ws.Range("BX47").Value = -0.00317995724487903
$ws.Range("BY47").Value = 0.4620566743634537
$ws.Range("BZ47").Value = 0.05304889653688862
$ws.Range("CA47").Value = -0.5170399564558569
$ws.Range("CB47").Value = 0.3901853270248822
$ws.Range("CC47").Value = -0.00217839596067916
$ws.Range("CD47").Value = 0.4496229560030233
$ws.Range("CE47").Value = -0.1136298432228148
$ws.Range("CF47").Value = -0.02895197187989974
$ws.Range("CG47").Value = 0.233384201608121
$ws.Range("CH47").Value = 0.0051319348870346
$ws.Range("CI47").Value = 0.0923147305357789
$ws.Range("CJ47").Value = 2.114242865113391
$ws.Range("CK47").Value = 0.0891408181828677
$ws.Range("CL47").Value = 7.925210431350729
$ws.Range("CM47").Value = -0.173274094326726
$ws.Range("CN47").Value = 12.69360902255639

# Row 65
$ws.Range("C65").Value = -0.1420493707021763
$ws.Range("D65").Value = -0.2064221155762598
$ws.Range("E65").Value = 0.3893617274955246
$ws.Range("F65").Value = 0.05141594900854674
$ws.Range("G65").Value = 1.314149824877073
$ws.Range("H65").Value = 0.3529792149415158
$ws.Range("I65").Value = 0.1486213121882455
$ws.Range("J65").Value = -0.5592112015157437
$ws.Range("K65").Value = 0.257659057139169
$ws.Range("L65").Value = 0.0786659134918908
$ws.Range("M65").Value = 0.2847943647262451
$ws.Range("N65").Value = 0.4666495748279185
$ws.Range("O65").Value = -0.07866590824398238
$ws.Range("P65").Value = 0.2847944916411569
$ws.Range("Q65").Value = 0.5333504251750465
$ws.Range("R65").Value = 85.33211705189514
$ws.Range("S65").Value = 0.2180548716611285
$ws.Range("T65").Value = 356.7282412170728
$ws.Range("U65").Value = -0.849583643994927
$ws.Range("V65").Value = 379.5840592510853
$ws.Range("Z65").Value = 1.046101686127638
$ws.Range("AA65").Value = 30.66398270595238
$ws.Range("AE65").Value = 1.056640691629629
$ws.Range("AF65").Value = 36.69698522058502
$ws.Range("AJ65").Value = 1.558086573414366
$ws.Range("AK65").Value = 30.61862178478972
$ws.Range("AO65").Value = 0.00430583118389643
$ws.Range("AP65").Value = 0.1250315298917528
$ws.Range("AT65").Value = 0.0043492105026945
$ws.Range("AU65").Value = 0.1496309285243018
$ws.Range("AY65").Value = 0.006413198491106671
$ws.Range("AZ65").Value = 0.1248465720072976
$ws.Range("BA65").Value = -1.740760819216694
$ws.Range("BB65").Value = -0.3063060096543673
$ws.Range("BC65").Value = 14.0087980266542
$ws.Range("BD65").Value = -0.02519534972365178
$ws.Range("BE65").Value = 14.68619528619529
$ws.Range("BF65").Value = -0.4712008016794442
$ws.Range("BG65").Value = -0.9948644078573503
$ws.Range("BH65").Value = 5.363052053236021
$ws.Range("BI65").Value = -0.0006098723079855535
$ws.Range("BJ65").Value = 5.386531986531986
$ws.Range("BK65").Value = -2.240183863204612
$ws.Range("BL65").Value = -0.3708365574906378
$ws.Range("BM65").Value = 6.02200311185707
$ws.Range("BN65").Value = 0.04002287021154942
$ws.Range("BO65").Value = 4.976430976430976
$ws.Range("BP65").Value = 0.9665526955008881
$ws.Range("BQ65").Value = -0.7834373325143476
$ws.Range("BR65").Value = 2.623742860256491
$ws.Range("BS65").Value = -0.0646083476272156
$ws.Range("BT65").Value = 4.323232323232324
$ws.Range("BU65").Value = 0.04845561027024558
$ws.Range("BV65").Value = 0.4492345122308752
$ws.Range("BW65").Value = 0.3653223231161853
$ws.Range("BX65").Value = 0.0005368789695481436
$ws.Range("BY65").Value = 0.3514847177321767
$ws.Range("BZ65").Value = -0.1378078235456782
$ws.Range("CA65").Value = -0.3678050298685199
$ws.Range("CB65").Value = 0.4231335885403462
$ws.Range("CC65").Value = 0.005152687661622438
$ws.Range("CD65").Value = 0.2875049821671612
$ws.Range("CE65").Value = 0.1104487407163886
$ws.Range("CF65").Value = -0.6934054006232365
$ws.Range("CG65").Value = 0.2115440882247533
$ws.Range("CH65").Value = -0.005689566631170593
$ws.Range("CI65").Value = 0.3610103001006619
$ws.Range("CJ65").Value = -1.679980682377926
$ws.Range("CK65").Value = 0.1216352577311278
$ws.Range("CL65").Value = 6.670441506703358
$ws.Range("CM65").Value = 0.004116638078902172
$ws.Range("CN65").Value = 6.557575757575757
